# Applies the commit "Added ABP test cases and modified IWP Bootstrap deferred test cases"
# to KatalonData/MultibillTestData/MultibillCCData.xlsx
#
# Updates the B2 (and B3, where present) "Date" timestamp cells on a number of
# per-scenario worksheets to reflect the latest test execution run, and flips
# the A2 "Result" cell from Fail -> Pass on the sheets whose deferred
# IWP Bootstrap test cases now pass.

$wb = $excel.ActiveWorkbook

function Set-DateCell {
    param($SheetName, $CellRef, $NewDate)
    $ws = $wb.Worksheets($SheetName)
    $ws.Range($CellRef).Value = $NewDate
}

function Set-ResultAndDate {
    param($SheetName, $NewResult, $NewDate)
    $ws = $wb.Worksheets($SheetName)
    $ws.Range("A2").Value = $NewResult
    $ws.Range("B2").Value = $NewDate
}

# Sheets whose Result is unchanged (Pass) - only the run Date is refreshed.
Set-DateCell "VerifyCANSearch"             "B2" "Mon Sep 15 21:12:28 IST 2025"
Set-DateCell "VerifyStreetAddressSearch"   "B2" "Mon Sep 15 21:13:44 IST 2025"
Set-DateCell "VerifyStaticTextOnViewCart"  "B2" "Mon Sep 15 21:25:18 IST 2025"
Set-DateCell "VerifyStaticTextOnViewCart"  "B3" "Mon Sep 15 21:26:00 IST 2025"
Set-DateCell "VerifyLookup1Search"         "B2" "Mon Sep 15 21:13:03 IST 2025"
Set-DateCell "VerifyUDF3Saerch"            "B2" "Mon Sep 15 21:14:26 IST 2025"
Set-DateCell "VerifyDataOnCartContent"     "B2" "Mon Sep 15 21:11:43 IST 2025"
Set-DateCell "VerifySearchResult"          "B2" "Mon Sep 15 21:24:37 IST 2025"
Set-DateCell "VerifyAmountTextBoxEditable" "B2" "Mon Sep 15 21:10:35 IST 2025"
Set-DateCell "VerifyALTIDRedacted"         "B2" "Mon Sep 15 21:09:01 IST 2025"
Set-DateCell "VerifyALTIDNotRedacted"      "B2" "Mon Sep 15 21:08:04 IST 2025"
Set-DateCell "VerifyStaticTextOnSearch"    "B2" "Mon Sep 15 21:26:42 IST 2025"
Set-DateCell "Verify2Pages"                "B2" "Mon Sep 15 21:27:13 IST 2025"

# Sheets whose deferred IWP Bootstrap test cases now pass: Result flips Fail -> Pass
# (and the run Date is refreshed at the same time).
Set-ResultAndDate "VerifyRemoveCartContent"    "Pass" "Mon Sep 15 21:23:35 IST 2025"
Set-ResultAndDate "VerifyPaymentEntryPageCC"   "Pass" "Mon Sep 15 21:18:39 IST 2025"
Set-ResultAndDate "VerifyPaymentEntryPagePC"   "Pass" "Mon Sep 15 21:22:29 IST 2025"
Set-ResultAndDate "VerifyPaymentEntryPageCorp" "Pass" "Mon Sep 15 21:19:47 IST 2025"
